$wb = $excel.ActiveWorkbook

# "settings" sheet: update the child form title string and make it the active sheet/selection
$settings = $wb.Worksheets.Item("settings")
$settings.Range("A2").Value = "Child Form V"

# Make "settings" the active sheet (moves tabSelected + activeTab there, removes it from "survey")
$settings.Activate()
$settings.Range("B2").Select()
